# Refresh cryptocurrency Price (D) and Volume(1h) (E) columns with the
# latest scraped values, per the scheduled GitHub Actions update run.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = "'88.178.99"
$ws.Range("E2").Value = '  -0.56%  '

$ws.Range("D3").Value = "'3.255.00"
$ws.Range("E3").Value = '  -3.71%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = "'212.55"
$ws.Range("E5").Value = '  -4.61%  '

$ws.Range("D6").Value = "'627.67"
$ws.Range("E6").Value = '  -5.43%  '

$ws.Range("D7").Value = "'0.385"
$ws.Range("E7").Value = '  +10.35%  '

$ws.Range("E8").Value = '  +14.62%  '

$ws.Range("E9").Value = '  +0.02%  '

$ws.Range("D10").Value = "'3.251.35"
$ws.Range("E10").Value = '  -3.77%  '

$ws.Range("D11").Value = "'0.573"
$ws.Range("E11").Value = '  -6.77%  '

$ws.Range("E12").Value = '  +12.19%  '

$ws.Range("E13").Value = '  -3.78%  '

$ws.Range("E14").Value = '  -0.38%  '

$ws.Range("D15").Value = "'34.14"
$ws.Range("E15").Value = '  -4.29%  '

$ws.Range("D16").Value = "'3.849.02"
$ws.Range("E16").Value = '  -3.54%  '

$ws.Range("D17").Value = "'87.854.57"
$ws.Range("E17").Value = '  -0.48%  '

$ws.Range("D18").Value = "'3.278.77"
$ws.Range("E18").Value = '  -2.59%  '

$ws.Range("D19").Value = "'3.24"
$ws.Range("E19").Value = '  -0.56%  '

$ws.Range("D20").Value = "'14.01"
$ws.Range("E20").Value = '  -6.00%  '

$ws.Range("D21").Value = "'434.95"
$ws.Range("E21").Value = '  -8.10%  '

$ws.Range("D22").Value = "'8.99"
$ws.Range("E22").Value = '  -3.91%  '

$ws.Range("D23").Value = "'5.35"
$ws.Range("E23").Value = '  -7.06%  '

$ws.Range("D24").Value = "'7.40"
$ws.Range("E24").Value = '  -2.25%  '

$ws.Range("D25").Value = "'5.38"
$ws.Range("E25").Value = '  -3.46%  '

$ws.Range("E26").Value = '  -9.16%  '

$ws.Range("D27").Value = "'0.0000143"
$ws.Range("E27").Value = '  +9.44%  '

$ws.Range("D29").Value = "'77.29"
$ws.Range("E29").Value = '  -3.41%  '

$ws.Range("E30").Value = '  -0.03%  '

$ws.Range("E31").Value = '  -13.95%  '

$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = '  +0.29%  '

$ws.Range("D33").Value = "'571.04"
$ws.Range("E33").Value = '  -6.31%  '

$ws.Range("D34").Value = "'8.86"
$ws.Range("E34").Value = '  -6.49%  '

$ws.Range("D35").Value = "'1.39"
$ws.Range("E35").Value = '  -12.33%  '

$ws.Range("D36").Value = "'7.25"
$ws.Range("E36").Value = '  +4.11%  '

$ws.Range("E37").Value = '  -6.56%  '

$ws.Range("E38").Value = '  -9.15%  '

$ws.Range("D39").Value = "'22.90"
$ws.Range("E39").Value = '  -5.89%  '

$ws.Range("D40").Value = "'3.30"
$ws.Range("E40").Value = '  +5.11%  '

$ws.Range("D41").Value = "'21.82"
$ws.Range("E41").Value = '  +0.58%  '

$ws.Range("E42").Value = '  -0.09%  '

$ws.Range("E43").Value = '  -6.23%  '

$ws.Range("E44").Value = '  -7.56%  '

$ws.Range("E45").Value = '  +0.00%  '

$ws.Range("D46").Value = "'151.54"
$ws.Range("E46").Value = '  -4.13%  '

$ws.Range("E47").Value = '  +18.13%  '

$ws.Range("D48").Value = "'180.11"
$ws.Range("E48").Value = '  -7.50%  '

$ws.Range("D49").Value = "'45.23"
$ws.Range("E49").Value = '  -5.94%  '

$ws.Range("E50").Value = '  -4.17%  '

$ws.Range("D51").Value = "'4.25"
$ws.Range("E51").Value = '  -3.61%  '
